$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cardholder name / account number / last name
$ws.Range("C2").Value = "Hartmut"

# Ensure card number stays text (it is a long digit string stored as text)
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "2570314725427075"
$ws.Range("C3").Value = "Mohaupt"

# Opening balance line
$ws.Range("D5").Value = "KONTOSTAND AM 18.01.2024"

# Row 6
$ws.Range("B6").Value = "19.01."
$ws.Range("C6").Value = "20.01."
$ws.Range("D6").Value = "ABSCHLAG STROM Stadtwerke Rosenheim 32937694"
$ws.Range("E6").Value = "85,94-"

# Row 7
$ws.Range("B7").Value = "21.01."
$ws.Range("C7").Value = "22.01."
$ws.Range("D7").Value = "BEITRAG Allianz SE K-49287274"
$ws.Range("E7").Value = "53,85-"

# Row 8
$ws.Range("B8").Value = "23.01."
$ws.Range("C8").Value = "24.01."
$ws.Range("D8").Value = "MITGLIEDSBEITRAG ZEUS BODYPOWER"
$ws.Range("E8").Value = "25,02-"

# Row 9
$ws.Range("B9").Value = "27.01."
$ws.Range("C9").Value = "28.01."
$ws.Range("D9").Value = "KARTENZAHLUNG JET TANKSTELLE"
$ws.Range("E9").Value = "45,35-"

# Closing balance line
$ws.Range("D12").Value = "KONTOSTAND AM 31.01.2024"
$ws.Range("E12").Value = "210,16-"

# Next billing date
$ws.Range("C13").Value = "IHR NAECHSTER ABRECHNUNGSTERMIN 05.02.2024"
